# Insert a new record row at row 295 (shifts existing rows 295-367 down to 296-368)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(295).Insert()

$ws.Range("A295").Value = 3
$ws.Range("B295").Value = "Femacal de La Calera"
$ws.Range("C295").Value = "Coquimbo"
$ws.Range("D295").Value = 44754
$ws.Range("E295").Value = 5
$ws.Range("F295").Value = 100112012
$ws.Range("G295").Value = "Espinaca"
$ws.Range("H295").Value = "Sin especificar"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 170
$ws.Range("K295").Value = 4000
$ws.Range("L295").Value = 4500
$ws.Range("M295").Value = 4235
$ws.Range("N295").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O295").Value = "Provincia de Quillota"
$ws.Range("P295").Value = 1412
$ws.Range("Q295").Value = 3
$ws.Range("R295").Value = "Hortaliza"
